# Update the "concise marksheet" correction/total marks figures on the
# quiz sheet, per commit message: "changes in concise marksheet / Corr/total marks"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking" -> B11 changes from 3 to 5
$ws.Range("B11").Value = 5

# Row 12 "Total" -> B12 changes from 57 to 95, and the Corr/Total text in E12
# changes from "53/84" to "95/140"
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
